$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06374189416929597
$ws.Range("H2").Value = 91.44597577355383
$ws.Range("I2").Value = 27.79781136887867
$ws.Range("G3").Value = 0.07280974543165851
$ws.Range("H3").Value = 44.27857189339247
$ws.Range("G4").Value = 0.03001691309245777
$ws.Range("H4").Value = -38.33301242710296
$ws.Range("G5").Value = 0.06997522669065702
$ws.Range("H5").Value = -0.9614648100257756
$ws.Range("G6").Value = -0.1265854152322287
$ws.Range("H6").Value = -7.617160112889162
$ws.Range("G7").Value = -0.137181892050301
$ws.Range("H7").Value = -9.716135762806777
$ws.Range("G8").Value = -0.2037904830513343
$ws.Range("H8").Value = -2.315387611540844
$ws.Range("G9").Value = -0.2911769043957947
$ws.Range("H9").Value = 3.608347163627863
$ws.Range("G10").Value = -0.006626872574338505
$ws.Range("H10").Value = -620.6878510981782
$ws.Range("G11").Value = 0.02257092478619185
$ws.Range("H11").Value = 209.762693887793
$ws.Range("G12").Value = 0.19910325564752
$ws.Range("H12").Value = -6.044049268094411
$ws.Range("G13").Value = 0.2531862596566299
$ws.Range("H13").Value = 7.740397556510701
$ws.Range("G14").Value = -0.06462633911388026
$ws.Range("H14").Value = 29.02125431619802
$ws.Range("G15").Value = -0.06081691571111515
$ws.Range("H15").Value = 14.30463265477377
$ws.Range("G16").Value = 0.1948543360970333
$ws.Range("H16").Value = 1.776442380774853
$ws.Range("G17").Value = 0.2139304205239499
$ws.Range("H17").Value = 23.04071088770731
$ws.Range("G18").Value = 0.05353976287189038
$ws.Range("H18").Value = -1.272534046304775
$ws.Range("G19").Value = 0.05366568516383025
$ws.Range("H19").Value = -37.62676214631674
$ws.Range("G20").Value = 0.01034344447897358
$ws.Range("H20").Value = -18.74051789024545
$ws.Range("G21").Value = -0.01760843828850251
$ws.Range("H21").Value = 67.31538936123839
$ws.Range("G22").Value = 0.06802829854133297
$ws.Range("H22").Value = 4.221606878597715
$ws.Range("G23").Value = 0.1186771471676996
$ws.Range("H23").Value = 105.7745334316112
$ws.Range("G24").Value = 0.03533730214771737
$ws.Range("H24").Value = 9.085577116980954
$ws.Range("G25").Value = 0.01590113315853392
$ws.Range("H25").Value = -45.97355696984602
$ws.Range("G26").Value = 0.09499952633269493
$ws.Range("H26").Value = -16.14587309076303
$ws.Range("G27").Value = 0.1070567709475595
$ws.Range("H27").Value = 18.70445907608494
$ws.Range("G28").Value = 0.1223076431759911
$ws.Range("H28").Value = 4.120750086914041
$ws.Range("G29").Value = 0.1454121878080378
$ws.Range("H29").Value = 21.55398034750331
$ws.Range("G30").Value = 0.04969523207167355
$ws.Range("H30").Value = -26.07992875163797
$ws.Range("G31").Value = 0.06038010479839933
$ws.Range("H31").Value = -12.00501010827067
$ws.Range("G32").Value = 0.05240961544770578
$ws.Range("H32").Value = 20.02232439871364
$ws.Range("G33").Value = 0.08700169883816945
$ws.Range("H33").Value = 60.11086548093487
$ws.Range("G34").Value = 0.007607627652831255
$ws.Range("H34").Value = 139.8333926129861
$ws.Range("G35").Value = 0.02443481770390726
$ws.Range("H35").Value = 75.12648272552137
$ws.Range("G36").Value = 0.004507954462800914
$ws.Range("H36").Value = -70.83913231169331
$ws.Range("G37").Value = 0.02173823189845242
$ws.Range("H37").Value = 73.57714912960435
$ws.Range("G38").Value = 0.06524856201101945
$ws.Range("H38").Value = -9.052110763348866
$ws.Range("G39").Value = 0.02966857644405549
$ws.Range("H39").Value = -31.10346044280285
$ws.Range("G40").Value = 0.06620111204117042
$ws.Range("H40").Value = 48.00236445347575
$ws.Range("G41").Value = 0.03943011885955505
$ws.Range("H41").Value = 219.0314007031886
$ws.Range("G42").Value = 0.07512414691143346
$ws.Range("H42").Value = 43.69857498301871
$ws.Range("G43").Value = 0.060577889450595
$ws.Range("H43").Value = 21.39885512877445
$ws.Range("G44").Value = 0.09472364758580461
$ws.Range("H44").Value = -28.10396586847564
$ws.Range("G45").Value = 0.1308004021031449
$ws.Range("H45").Value = -27.10092081164279
$ws.Range("G46").Value = -0.02508464513799502
$ws.Range("H46").Value = 42.89839439327003
$ws.Range("G47").Value = 0.02691992482108009
$ws.Range("H47").Value = 1127.647965338288
$ws.Range("G48").Value = -0.0001689302114396751
$ws.Range("H48").Value = -101.1656934922414
$ws.Range("G49").Value = 0.00007893714018425944
$ws.Range("H49").Value = 101.4200345104302
$ws.Range("G50").Value = 0.1430713303841797
$ws.Range("H50").Value = 0.09785091779785221
$ws.Range("G51").Value = 0.1361590792330044
$ws.Range("H51").Value = 3.965835247552169
$ws.Range("G52").Value = 0.06746770057011559
$ws.Range("H52").Value = 8.901945954531479
$ws.Range("G53").Value = 0.06563826013877397
$ws.Range("H53").Value = 7.28426466587395
$ws.Range("G54").Value = -0.1312649730868879
$ws.Range("H54").Value = -46.98182840663135
$ws.Range("G55").Value = -0.04908029145359365
$ws.Range("H55").Value = 52.71038214676504
$ws.Range("G56").Value = 0.1345192637838271
$ws.Range("H56").Value = -13.28254219933515
$ws.Range("G57").Value = 0.1762308217490755
$ws.Range("H57").Value = 26.34292944341037
